$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($ws, $addr, $val)
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextCell $ws "C2" "21.0"
Set-TextCell $ws "E2" "Lugar poblado"
Set-TextCell $ws "F2" "2.3"
Set-TextCell $ws "G2" "9.906393"
Set-TextCell $ws "H2" "-83.67443"
Set-TextCell $ws "I2" "CARTAGO"
Set-TextCell $ws "J2" "Peligroso"
Set-TextCell $ws "A3" "03/06/2020"
Set-TextCell $ws "B3" "23:55:17"
Set-TextCell $ws "C3" "21.0"
Set-TextCell $ws "D3" "SUBDUCCION_PLACA"
Set-TextCell $ws "E3" "Leve en Jacó"
Set-TextCell $ws "F3" "3.4"
Set-TextCell $ws "G3" "9.631"
Set-TextCell $ws "H3" "-84.7381"
Set-TextCell $ws "I3" "PUNTARENAS"
Set-TextCell $ws "J3" "11.7 km oeste de Jacó,"
Set-TextCell $ws "A4" "01/06/2020"
Set-TextCell $ws "B4" "19:34:00"
Set-TextCell $ws "C4" "10.0"
Set-TextCell $ws "D4" "TECTONICO_POR_FALLA_LOCAL"
Set-TextCell $ws "E4" "Leve en Turrialba centro"
Set-TextCell $ws "F4" "2.2"
Set-TextCell $ws "G4" "9.981"
Set-TextCell $ws "H4" "-83.71"
Set-TextCell $ws "I4" "CARTAGO"
Set-TextCell $ws "J4" "8 km NO"
Set-TextCell $ws "A5" "26/04/2020"
Set-TextCell $ws "B5" "00:14:15"
Set-TextCell $ws "C5" "21.0"
Set-TextCell $ws "E5" "Reportado en: Puerto Carrillo, Punta Islita, Playa Sámara, Huacas, Hojancha, Nosara, Nicoya, Nandayure, Nambí, Filadelfia, Santa Cruz"
Set-TextCell $ws "F5" "4.4"
Set-TextCell $ws "G5" "9.8414"
Set-TextCell $ws "H5" "-85.5104"
Set-TextCell $ws "I5" "GUANACASTE"
Set-TextCell $ws "J5" "4.4 km Suroeste de Puerto Carrillo"
Set-TextCell $ws "A6" "20/04/2020"
Set-TextCell $ws "B6" "03:27:00"
Set-TextCell $ws "C6" "28.0"
Set-TextCell $ws "D6" "DEFORMACION_INTERNA"
Set-TextCell $ws "E6" "Reportado en :Talamanca, Valle de La Estrella y Limón"
Set-TextCell $ws "F6" "4.0"
Set-TextCell $ws "G6" "9.71"
Set-TextCell $ws "H6" "-82.8487"
Set-TextCell $ws "I6" "LIMON"
Set-TextCell $ws "J6" "3 km al Suroeste de Cahuita"
Set-TextCell $ws "A7" "30/12/2017"
Set-TextCell $ws "B7" "07:38:00"
Set-TextCell $ws "C7" "10.0"
Set-TextCell $ws "D7" "TECTONICO_POR_FALLA_LOCAL"
Set-TextCell $ws "E7" "detalle"
Set-TextCell $ws "F7" "4.4"
Set-TextCell $ws "G7" "9.33"
Set-TextCell $ws "H7" "-83.2865"
Set-TextCell $ws "I7" "PUNTARENAS"
Set-TextCell $ws "J7" "17.4KM Noroeste"
Set-TextCell $ws "A8" "30/12/2017"
Set-TextCell $ws "B8" "17:37:00"
Set-TextCell $ws "C8" "9.0"
Set-TextCell $ws "D8" "TECTONICO_POR_FALLA_LOCAL"
Set-TextCell $ws "F8" "3.6"
Set-TextCell $ws "G8" "9.68"
Set-TextCell $ws "H8" "-84.87"
Set-TextCell $ws "I8" "PUNTARENAS"
Set-TextCell $ws "J8" "17 Km Sede vaquera"
Set-TextCell $ws "A9" "28/12/2017"
Set-TextCell $ws "B9" "17:01:00"
Set-TextCell $ws "C9" "147.0"
Set-TextCell $ws "E9" "reportado en atena"
Set-TextCell $ws "F9" "4.6"
Set-TextCell $ws "G9" "10.6358"
Set-TextCell $ws "H9" "-84.7997"
Set-TextCell $ws "I9" "SAN_JOSE"
Set-TextCell $ws "J9" "15.5 km noroeste de Tilaran"
Set-TextCell $ws "A10" "26/12/2017"
Set-TextCell $ws "B10" "17:20:00"
Set-TextCell $ws "C10" "10.0"
Set-TextCell $ws "F10" "3.9"
Set-TextCell $ws "G10" "9.46"
Set-TextCell $ws "H10" "-84.59"
Set-TextCell $ws "J10" "sureste de jaco"
Set-TextCell $ws "I12" "HEREDIA"
Set-TextCell $ws "I16" "ALAJUELA"
Set-TextCell $ws "I22" "HEREDIA"
Set-TextCell $ws "I23" "HEREDIA"
Set-TextCell $ws "I29" "HEREDIA"
Set-TextCell $ws "I31" "HEREDIA"
Set-TextCell $ws "A35" "24/02/2012"
Set-TextCell $ws "B35" "23:07:46"
Set-TextCell $ws "C35" "8.7"
Set-TextCell $ws "D35" "CHOQUE_PLACAS"
Set-TextCell $ws "E35" "detalle innecesario"
Set-TextCell $ws "F35" "2.7"
Set-TextCell $ws "G35" "4.9"
Set-TextCell $ws "H35" "67.0"
Set-TextCell $ws "I35" "PUNTARENAS"
Set-TextCell $ws "J35" "Descripcion"
Set-TextCell $ws "A36" "24/12/2014"
Set-TextCell $ws "B36" "23:07:46"
Set-TextCell $ws "C36" "8.7"
Set-TextCell $ws "D36" "TECTONICO_SUBDUCCION"
Set-TextCell $ws "E36" "dt"
Set-TextCell $ws "F36" "5.0"
Set-TextCell $ws "G36" "4.9"
Set-TextCell $ws "H36" "67.0"
Set-TextCell $ws "I36" "ALAJUELA"
Set-TextCell $ws "J36" "Descripcion"
Set-TextCell $ws "A37" "19/10/2019"
Set-TextCell $ws "B37" "07:06:56"
Set-TextCell $ws "C37" "45.7"
Set-TextCell $ws "D37" "DEFORMACION_INTERNA"
Set-TextCell $ws "E37" "Isla de Chira"
Set-TextCell $ws "F37" "4.41"
Set-TextCell $ws "G37" "10.0974"
Set-TextCell $ws "H37" "-85.1278"
Set-TextCell $ws "I37" "SIN_ASIGNAR"
Set-TextCell $ws "J37" "Marino"
